$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "target column" data row (row 62): B62 first so its
# (unaligned) style is minted before A62's (left-aligned) style -
# keeps cellXfs insertion order matching the authored file.
$ws.Range("B62").Value = 0
$ws.Range("B62").Font.Name = "Arial"
$ws.Range("B62").Font.Size = 11

$ws.Range("A62").Value = "131306-0.0"
$ws.Range("A62").Font.Name = "Arial"
$ws.Range("A62").Font.Size = 11
$ws.Range("A62").HorizontalAlignment = -4131

# Reflect the new selection / scrolled view position used when the row
# was added.
$ws.Activate()
$excel.Goto($ws.Range("A57"), $true) | Out-Null
$ws.Range("A62").Select() | Out-Null
